$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 28.48226033333333
$ws.Range("H2").Value = 85.44678099999999
$ws.Range("I2").Value = 0.2101651977164657
$ws.Range("J2").Value = 0.2101651977164658
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.2423163333333333
$ws.Range("N2").Value = 0.726949
$ws.Range("O2").Value = 0.08545736829972225
$ws.Range("P2").Value = 0.08545736829972224
$ws.Range("Q2").Value = 6.901716889018776
$ws.Range("R2").Value = 62.11545200116899
$ws.Range("S2").Value = 0.01796016470503996
$ws.Range("T2").Value = 0.01796016470503996
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 28.48226033333333
$ws.Range("H3").Value = 85.44678099999999
$ws.Range("I3").Value = 0.2101651977164657
$ws.Range("J3").Value = 0.2101651977164658
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.274826333333333
$ws.Range("N3").Value = 6.824479
$ws.Range("O3").Value = 0.8022598770432592
$ws.Range("P3").Value = 0.8022598770432591
$ws.Range("Q3").Value = 64.7921958391221
$ws.Range("R3").Value = 583.1297625520989
$ws.Range("S3").Value = 0.1686071056787841
$ws.Range("T3").Value = 0.1686071056787841
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 28.48226033333333
$ws.Range("H4").Value = 85.44678099999999
$ws.Range("I4").Value = 0.2101651977164657
$ws.Range("J4").Value = 0.2101651977164658
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3183803333333333
$ws.Range("N4").Value = 0.955141
$ws.Range("O4").Value = 0.1122827546570186
$ws.Range("P4").Value = 0.1122827546570186
$ws.Range("Q4").Value = 9.068191539013442
$ws.Range("R4").Value = 81.61372385112099
$ws.Range("S4").Value = 0.02359792733264173
$ws.Range("T4").Value = 0.02359792733264173
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 83.45109033333334
$ws.Range("H5").Value = 250.353271
$ws.Range("I5").Value = 0.6157697701763504
$ws.Range("J5").Value = 0.6157697701763504
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.2423163333333333
$ws.Range("N5").Value = 0.726949
$ws.Range("O5").Value = 0.08545736829972225
$ws.Range("P5").Value = 0.08545736829972224
$ws.Range("Q5").Value = 20.22156222224211
$ws.Range("R5").Value = 181.994060000179
$ws.Range("S5").Value = 0.0526220640377957
$ws.Range("T5").Value = 0.05262206403779569
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 83.45109033333334
$ws.Range("H6").Value = 250.353271
$ws.Range("I6").Value = 0.6157697701763504
$ws.Range("J6").Value = 0.6157697701763504
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.274826333333333
$ws.Range("N6").Value = 6.824479
$ws.Range("O6").Value = 0.8022598770432592
$ws.Range("P6").Value = 0.8022598770432591
$ws.Range("Q6").Value = 189.8367378356455
$ws.Range("R6").Value = 1708.530640520809
$ws.Range("S6").Value = 0.4940073801086348
$ws.Range("T6").Value = 0.4940073801086348
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 83.45109033333334
$ws.Range("H7").Value = 250.353271
$ws.Range("I7").Value = 0.6157697701763504
$ws.Range("J7").Value = 0.6157697701763504
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.3183803333333333
$ws.Range("N7").Value = 0.955141
$ws.Range("O7").Value = 0.1122827546570186
$ws.Range("P7").Value = 0.1122827546570186
$ws.Range("Q7").Value = 26.56918595735678
$ws.Range("R7").Value = 239.122673616211
$ws.Range("S7").Value = 0.06914032602991987
$ws.Range("T7").Value = 0.06914032602991987
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.78130366666667
$ws.Range("H8").Value = 44.343911
$ws.Range("I8").Value = 0.1090684366779874
$ws.Range("J8").Value = 0.1090684366779875
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2423163333333333
$ws.Range("N8").Value = 0.726949
$ws.Range("O8").Value = 0.08545736829972225
$ws.Range("P8").Value = 0.08545736829972224
$ws.Range("Q8").Value = 3.581751306393222
$ws.Range("R8").Value = 32.235761757539
$ws.Range("S8").Value = 0.009320701563065708
$ws.Range("T8").Value = 0.009320701563065708
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.78130366666667
$ws.Range("H9").Value = 44.343911
$ws.Range("I9").Value = 0.1090684366779874
$ws.Range("J9").Value = 0.1090684366779875
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.274826333333333
$ws.Range("N9").Value = 6.824479
$ws.Range("O9").Value = 0.8022598770432592
$ws.Range("P9").Value = 0.8022598770432591
$ws.Range("Q9").Value = 33.62489882192989
$ws.Range("R9").Value = 302.624089397369
$ws.Range("S9").Value = 0.08750123059858271
$ws.Range("T9").Value = 0.08750123059858271
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.78130366666667
$ws.Range("H10").Value = 44.343911
$ws.Range("I10").Value = 0.1090684366779874
$ws.Range("J10").Value = 0.1090684366779875
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.3183803333333333
$ws.Range("N10").Value = 0.955141
$ws.Range("O10").Value = 0.1122827546570186
$ws.Range("P10").Value = 0.1122827546570186
$ws.Range("Q10").Value = 4.706076388494555
$ws.Range("R10").Value = 42.354687496451
$ws.Range("S10").Value = 0.01224650451633903
$ws.Range("T10").Value = 0.01224650451633903
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 8.808546666666667
$ws.Range("H11").Value = 26.42564
$ws.Range("I11").Value = 0.0649965954291964
$ws.Range("J11").Value = 0.06499659542919642
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.2423163333333333
$ws.Range("N11").Value = 0.726949
$ws.Range("O11").Value = 0.08545736829972225
$ws.Range("P11").Value = 0.08545736829972224
$ws.Range("Q11").Value = 2.134454730262222
$ws.Range("R11").Value = 19.21009257236
$ws.Range("S11").Value = 0.005554437993820881
$ws.Range("T11").Value = 0.005554437993820881
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 8.808546666666667
$ws.Range("H12").Value = 26.42564
$ws.Range("I12").Value = 0.0649965954291964
$ws.Range("J12").Value = 0.06499659542919642
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2.274826333333333
$ws.Range("N12").Value = 6.824479
$ws.Range("O12").Value = 0.8022598770432592
$ws.Range("P12").Value = 0.8022598770432591
$ws.Range("Q12").Value = 20.03791391572889
$ws.Range("R12").Value = 180.34122524156
$ws.Range("S12").Value = 0.05214416065725757
$ws.Range("T12").Value = 0.05214416065725758
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 8.808546666666667
$ws.Range("H13").Value = 26.42564
$ws.Range("I13").Value = 0.0649965954291964
$ws.Range("J13").Value = 0.06499659542919642
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.3183803333333333
$ws.Range("N13").Value = 0.955141
$ws.Range("O13").Value = 0.1122827546570186
$ws.Range("P13").Value = 0.1122827546570186
$ws.Range("Q13").Value = 2.804468023915555
$ws.Range("R13").Value = 25.24021221524
$ws.Range("S13").Value = 0.007297996778117956
$ws.Range("T13").Value = 0.007297996778117957
